{"js": "// Change the course-number heading from \"7.4 \" to \"10.4 \" (course moved\n// to section 10.4 for SP24), keeping the existing bold/size-24 run\n// formatting that the heading already uses.\nconst body = context.document.body;\n\nconst results = body.search(\"7.4 \", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find heading text \"7.4 \" to update.');\n}\n\n// Replace just the text in place; Office.js keeps the formatting of the\n// run(s) being replaced, so the bold / 24pt heading style is preserved.\nresults.items[0].insertText(\"10.4 \", Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Change the course-number heading from \"7.4 \" to \"10.4 \" (course moved\n# to section 10.4 for SP24). Find/Replace only touches the matched text,\n# so the heading's existing bold / 24pt run formatting is preserved.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$found = $find.Execute(\"7.4 \", $false, $false, $false, $false, $false, $true, 1, $false, \"10.4 \", 2)\n\nif (-not $found) {\n    throw 'Could not find heading text \"7.4 \" to update.'\n}\n"}
